$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 1.6
$ws.Range("I4").Value = 5.25
$ws.Range("S4").Value = 1.93
$ws.Range("T4").Value = 1.93
$ws.Range("W4").Value = 3.25
$ws.Range("X4").Value = 1.33
$ws.Range("AC4").Value = 6.5
$ws.Range("AH4").Value = 29
$ws.Range("AM4").Value = 351

# Row 5 updates
$ws.Range("G5").Value = 4.3
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 1.87
$ws.Range("J5").Value = 4.65
$ws.Range("K5").Value = 1.98
$ws.Range("L5").Value = 2.5
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 2.95
$ws.Range("S5").Value = 1.87
$ws.Range("T5").Value = 1.75
$ws.Range("W5").Value = 2.95
$ws.Range("X5").Value = 1.29
$ws.Range("Y5").Value = 1.44
$ws.Range("Z5").Value = 2.42
$ws.Range("AC5").Value = 12
$ws.Range("AD5").Value = 26
$ws.Range("AE5").Value = 13.5
$ws.Range("AF5").Value = 75
$ws.Range("AH5").Value = 40
$ws.Range("AI5").Value = 9
$ws.Range("AJ5").Value = 6.1
$ws.Range("AM5").Value = 450
$ws.Range("AN5").Value = 6.9
$ws.Range("AO5").Value = 9
$ws.Range("AQ5").Value = 16.5
$ws.Range("AR5").Value = 15
